# Update countries & provincias Spain
#
# The "Pais" sheet is a COVID-19 stats table (columns A..H = Pais,
# Casos totales, Nuevos casos, Casos activos, Recuperados, Casos
# criticos, Muertes hoy, Muertes) kept sorted by column B (Casos
# totales) descending. This refreshes a handful of countries' figures
# and the "last updated" banner; because a few countries' totals
# crossed each other, three rows (148-150) and two rows (214-215) swap
# country names to keep the list sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Pais")

# Banner timestamp (A1)
$ws.Range("A1").Value = 'Datos actualizados a 10 de Septiembre de 2020 a las 16:08'

# Estados Unidos (row 4)
$ws.Range("B4").Value = 6550984
$ws.Range("C4").Value = 1509
$ws.Range("E4").Value = 2509008

# Irak (row 23)
$ws.Range("B23").Value = 278418
$ws.Range("C23").Value = 4597
$ws.Range("D23").Value = 213817
$ws.Range("E23").Value = 56787
$ws.Range("G23").Value = 82
$ws.Range("H23").Value = 7814

# Alemania (row 24)
$ws.Range("B24").Value = 256835
$ws.Range("C24").Value = 486
$ws.Range("E24").Value = 15523
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 9412

# Emiratos Arabes Unidos (row 46)
$ws.Range("B46").Value = 76911
$ws.Range("C46").Value = 930
$ws.Range("D46").Value = 67945
$ws.Range("E46").Value = 8568
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 398

# Suiza (row 60)
$ws.Range("E60").Value = 5591
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 2020

# Serbia (row 69)
$ws.Range("B69").Value = 32136
$ws.Range("C69").Value = 58
$ws.Range("D69").Value = 30943
$ws.Range("E69").Value = 464
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 729

# Libia (row 80)
$ws.Range("B80").Value = 20939
$ws.Range("C80").Value = 477
$ws.Range("D80").Value = 2420
$ws.Range("E80").Value = 18180
$ws.Range("G80").Value = 15
$ws.Range("H80").Value = 339

# Noruega (row 92)
$ws.Range("B92").Value = 11801
$ws.Range("C92").Value = 55
$ws.Range("E92").Value = 2188
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 265

# Tayikistan (row 100)
$ws.Range("B100").Value = 8939
$ws.Range("C100").Value = 40
$ws.Range("D100").Value = 7714
$ws.Range("E100").Value = 1153
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 72

# Sri Lanka (row 133)
$ws.Range("B133").Value = 3152
$ws.Range("C133").Value = 5
$ws.Range("E133").Value = 185

# Rows 148-150 updated & re-sorted: Birmania overtakes Botsuana and
# Sierra Leona in total cases, so the three countries shift up one slot.
$ws.Range("A148").Value = 'Birmania'
$ws.Range("B148").Value = 2150
$ws.Range("C148").Value = 261
$ws.Range("D148").Value = 625
$ws.Range("E148").Value = 1511
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 14

$ws.Range("A149").Value = 'Botsuana'
$ws.Range("B149").Value = 2126
$ws.Range("D149").Value = 493
$ws.Range("E149").Value = 1624
$ws.Range("H149").Value = 9

$ws.Range("A150").Value = 'Sierra Leona'
$ws.Range("B150").Value = 2067
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 1622
$ws.Range("E150").Value = 373
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 72

# Liechtenstein (row 195)
$ws.Range("B195").Value = 109
$ws.Range("C195").Value = 1
$ws.Range("E195").Value = 3

# Rows 214-215: Montserrat overtakes Islas Malvinas, so they swap.
$ws.Range("A214").Value = 'Montserrat'
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = 'Islas Malvinas'
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
